$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 138, pushing the existing rows
# 138..161 down to 140..163 (formatting is inherited from the row above,
# which keeps column D's date style "s=2").
$ws.Rows("138:139").Insert()

# Row 138 - new weekly entry (Primera)
$ws.Range("A138").Value = 4
$ws.Range("B138").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C138").Value = "Los Lagos"
$ws.Range("D138").Value = 44504
$ws.Range("E138").Value = 10
$ws.Range("F138").Value = 100112021
$ws.Range("G138").Value = "Ají"
$ws.Range("H138").Value = "Inferno"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 40
$ws.Range("K138").Value = 33000
$ws.Range("L138").Value = 33000
$ws.Range("M138").Value = 33000
$ws.Range("N138").Value = "$/caja 12 kilos"
$ws.Range("O138").Value = "Región de Arica y Parinacota"
$ws.Range("P138").Value = 2750
$ws.Range("Q138").Value = 12
$ws.Range("R138").Value = "Hortaliza"

# Row 139 - new weekly entry (Segunda)
$ws.Range("A139").Value = 4
$ws.Range("B139").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C139").Value = "Los Lagos"
$ws.Range("D139").Value = 44504
$ws.Range("E139").Value = 10
$ws.Range("F139").Value = 100112021
$ws.Range("G139").Value = "Ají"
$ws.Range("H139").Value = "Inferno"
$ws.Range("I139").Value = "Segunda"
$ws.Range("J139").Value = 40
$ws.Range("K139").Value = 26000
$ws.Range("L139").Value = 26000
$ws.Range("M139").Value = 26000
$ws.Range("N139").Value = "$/caja 12 kilos"
$ws.Range("O139").Value = "Región de Arica y Parinacota"
$ws.Range("P139").Value = 2167
$ws.Range("Q139").Value = 12
$ws.Range("R139").Value = "Hortaliza"
